$d = $word.ActiveDocument

$replacements = @(
    @{old="870×8=6960"; new="321×6=1926"},
    @{old="780×8=6240"; new="922×3=2766"},
    @{old="548×2=1096"; new="693×2=1386"},
    @{old="780×6=4680"; new="313×2=626"},
    @{old="222×3=666";  new="709×4=2836"},
    @{old="469×2=938";  new="525×5=2625"},
    @{old="459×7=3213"; new="398×5=1990"},
    @{old="924×6=5544"; new="377×3=1131"},
    @{old="620×5=3100"; new="753×8=6024"},
    @{old="718×3=2154"; new="120×4=480"},
    @{old="102×7=714";  new="780×7=5460"},
    @{old="952×3=2856"; new="733×2=1466"},
    @{old="441×9=3969"; new="694×4=2776"},
    @{old="683×2=1366"; new="618×5=3090"},
    @{old="243×2=486";  new="882×5=4410"},
    @{old="173×2=346";  new="600×9=5400"},
    @{old="804×7=5628"; new="984×7=6888"},
    @{old="813×5=4065"; new="295×9=2655"},
    @{old="264×9=2376"; new="319×3=957"},
    @{old="316×8=2528"; new="405×4=1620"},
    @{old="708×2=1416"; new="836×3=2508"},
    @{old="479×2=958";  new="919×5=4595"},
    @{old="543×7=3801"; new="251×3=753"},
    @{old="972×5=4860"; new="308×9=2772"},
    @{old="208×6=1248"; new="204×9=1836"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
